$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68; this shifts the existing rows 68-82 down to 69-83
$ws.Rows.Item(68).EntireRow.Insert()

# Populate the newly inserted row 68 with the new record's data
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(68, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value = Get-Date -Year 2021 -Month 12 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(68, 5).Value = 15
$ws.Cells.Item(68, 6).Value = "Fruta"
$ws.Cells.Item(68, 7).Value = 100102
$ws.Cells.Item(68, 8).Value = "Cítricos"
$ws.Cells.Item(68, 9).Value = 100102004
$ws.Cells.Item(68, 10).Value = "Mandarina"
$ws.Cells.Item(68, 11).Value = "Murcott"
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 250
$ws.Cells.Item(68, 14).Value = 14000
$ws.Cells.Item(68, 15).Value = 15000
$ws.Cells.Item(68, 16).Value = 14500
$ws.Cells.Item(68, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(68, 18).Value = "Región Metropolitana"
$ws.Cells.Item(68, 19).Value = 725
$ws.Cells.Item(68, 20).Value = 20
